# Applies cryptos list update (prices/volumes refresh + two coin-row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "85.084.44"
$ws.Range("E2").Value = "  +5.29%  "
$ws.Range("D3").Value = "3.310.34"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'218.38"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").Value = "'636.00"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "'0.321"
$ws.Range("E7").Value = "  +11.57%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.592"
$ws.Range("E9").Value = "  -2.70%  "
$ws.Range("D10").Value = "3.301.99"
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("D11").Value = "'0.596"
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("D12").Value = "'0.0000275"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "3.909.44"
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("D15").Value = "'34.24"
$ws.Range("E15").Value = "  +3.94%  "
$ws.Range("D16").Value = "'5.41"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").Value = "84.760.23"
$ws.Range("E17").Value = "  +5.32%  "
$ws.Range("D18").Value = "3.297.95"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "'14.60"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").Value = "'3.16"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("D21").Value = "'439.85"
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("D22").Value = "'9.17"
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").Value = "'5.25"
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("D24").Value = "'7.38"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "'5.48"
$ws.Range("E25").Value = "  +12.16%  "
$ws.Range("D26").Value = "'12.18"
$ws.Range("E26").Value = "  +9.75%  "
$ws.Range("D27").Value = "3.473.45"
$ws.Range("E27").Value = "  +1.90%  "
$ws.Range("D28").Value = "'78.04"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "'604.86"
$ws.Range("E31").Value = "  +6.10%  "
$ws.Range("B32").Value = "Cronos"
$ws.Range("C32").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D32").Value = "'0.163"
$ws.Range("E32").Value = "  +30.68%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'9.21"
$ws.Range("E33").Value = "  -2.12%  "
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("D36").Value = "'2.04"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -6.22%  "
$ws.Range("D38").Value = "'23.19"
$ws.Range("E38").Value = "  -2.90%  "
$ws.Range("D39").Value = "'6.42"
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.418"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'20.98"
$ws.Range("E42").Value = "  +3.15%  "
$ws.Range("D43").Value = "'2.03"
$ws.Range("E43").Value = "  +9.51%  "
$ws.Range("D44").Value = "'3.05"
$ws.Range("E44").Value = "  +9.06%  "
$ws.Range("D45").Value = "'159.95"
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'189.29"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").Value = "'1.36"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").Value = "'44.87"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("D50").Value = "'0.786"
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").Value = "'26.43"
$ws.Range("E51").Value = "  +1.44%  "
